$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CategoryMap")

# Update individual cell values on existing rows
$ws.Range("G2").Value = 2
$ws.Range("G7").Value = 8
$ws.Range("G8").Value = 4
$ws.Range("G9").Value = 4
$ws.Range("G11").Value = 6
$ws.Range("G13").Value = 6
$ws.Range("D15").Value = 6

# Add a new (mostly empty) row 20, matching style of the existing data rows
$ws.Range("B20:J20").WrapText = $true

# Update the selection to match the post-edit state (row 7 selected)
$ws.Rows(7).Select()
